# The second table ("Äquivalenzklasse / Repräsentant / pass / fail") has
# a header row followed by six data rows. For each data row, mark the
# "pass" column (3rd column) with a centered "x" and center the
# (still empty) "fail" column (4th column) as well.

$d = $word.ActiveDocument
$t = $d.Tables(2)

for ($row = 2; $row -le 7; $row++) {
    $passCell = $t.Cell($row, 3)
    $passCell.Range.Text = "x"
    $passCell.Range.ParagraphFormat.Alignment = 1

    $failCell = $t.Cell($row, 4)
    $failCell.Range.ParagraphFormat.Alignment = 1
}
